$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated sensor readings (Time / Dust Level columns), row 3 now holds a
# mismatched Time/Dust-Level pair per the source export.
$updates = @(
    @{ Row = 2; Col = 1; Value = "20:5:27" },
    @{ Row = 2; Col = 2; Value = "22.804822310888497" },
    @{ Row = 3; Col = 1; Value = "0.585938" },
    @{ Row = 3; Col = 2; Value = "20:5:27" },
    @{ Row = 4; Col = 1; Value = "20:5:27" },
    @{ Row = 4; Col = 2; Value = "4.987632144992207" }
)

foreach ($u in $updates) {
    $cell = $ws.Cells.Item($u.Row, $u.Col)
    # Force text storage so numeric-looking / time-looking strings aren't
    # reinterpreted as numbers or date/time serials, then restore the
    # cell's normal (unformatted) style.
    $cell.NumberFormat = "@"
    $cell.Value2 = $u.Value
    $cell.Style = "Normal"
}
